# SectorGroup.xlsx: the codeforiati:category-name/group-name columns (D/E)
# and the codeforiati:group-code/category-code columns (F/G) were each
# other's values - this swaps D<->E and F<->G for every used row (including
# the header row) so the header/value pairing lines up correctly.
#
# Values in F/G are textual codes (e.g. "110", "111") stored as shared
# strings, not numbers, so plain Range.Value re-entry (which would coerce
# numeric-looking text to a Double) is avoided in favour of Copy/PasteSpecial,
# which moves the cell's stored type as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$rngD = $ws.Range("D1:D$lastRow")
$rngE = $ws.Range("E1:E$lastRow")
$rngF = $ws.Range("F1:F$lastRow")
$rngG = $ws.Range("G1:G$lastRow")
$rngScratch = $ws.Range("J1:J$lastRow")

# Swap column D and column E
$rngD.Copy()
$rngScratch.PasteSpecial()
$rngE.Copy()
$rngD.PasteSpecial()
$rngScratch.Copy()
$rngE.PasteSpecial()
$rngScratch.ClearContents()

# Swap column F and column G
$rngF.Copy()
$rngScratch.PasteSpecial()
$rngG.Copy()
$rngF.PasteSpecial()
$rngScratch.Copy()
$rngG.PasteSpecial()
$rngScratch.ClearContents()
